# Week 13 logging update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "YDS" - append this week's play-by-play yardage logs (space
# separated numbers stored as text) for OFF/DEF rush & pass.
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Range("B2").Value2 = $ydsWs.Range("B2").Value2 + " 4 2 4 2 4 0 5 7 0 7 5 5 1 7 16 2 2 8 3"
$ydsWs.Range("B3").Value2 = $ydsWs.Range("B3").Value2 + " 15 5 12 11 2 10 4 14 2 5"
$ydsWs.Range("C2").Value2 = $ydsWs.Range("C2").Value2 + " 2 10 8 6 2 1 10 1 -1 6 0 7 6 0 3 8 2 0 7 0 -1 15 11 -1 4 9 3 16 4 1 11 25 4 5 3 6 -3 7 7 3 6 10 0 12 3"
$ydsWs.Range("C3").Value2 = $ydsWs.Range("C3").Value2 + " 12 4 8 9 14 6 13 2 8 17 24 14 4 17 0 6"

# ---------------------------------------------------------------------
# Sheet "OFF" - cumulative offensive stat totals
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("B2").Value2 = 9
$offWs.Range("C2").Value2 = 257
$offWs.Range("E2").Value2 = 17
$offWs.Range("F2").Value2 = 111
$offWs.Range("G2").Value2 = 73
$offWs.Range("H2").Value2 = 7
$offWs.Range("I2").Value2 = 11
$offWs.Range("J2").Value2 = 53
$offWs.Range("L2").Value2 = 474
$offWs.Range("M2").Value2 = 325
$offWs.Range("O2").Value2 = 33
$offWs.Range("Q2").Value2 = 817

$offWs.Range("B3").Value2 = 16
$offWs.Range("C3").Value2 = 286
$offWs.Range("E3").Value2 = 52
$offWs.Range("F3").Value2 = 149
$offWs.Range("G3").Value2 = 50
$offWs.Range("H3").Value2 = 50
$offWs.Range("I3").Value2 = 95
$offWs.Range("J3").Value2 = 102
$offWs.Range("N3").Value2 = 43

# ---------------------------------------------------------------------
# Sheet "DEF" - cumulative defensive stat totals
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("B2").Value2 = 7
$defWs.Range("C2").Value2 = 390
$defWs.Range("D2").Value2 = 25
$defWs.Range("E2").Value2 = 15
$defWs.Range("F2").Value2 = 103
$defWs.Range("G2").Value2 = 111
$defWs.Range("I2").Value2 = 9
$defWs.Range("J2").Value2 = 55
$defWs.Range("L2").Value2 = 461
$defWs.Range("M2").Value2 = 312
$defWs.Range("O2").Value2 = 35
$defWs.Range("P2").Value2 = 20
$defWs.Range("Q2").Value2 = 909

$defWs.Range("B3").Value2 = 19
$defWs.Range("C3").Value2 = 277
$defWs.Range("D3").Value2 = 9
$defWs.Range("E3").Value2 = 49
$defWs.Range("F3").Value2 = 179
$defWs.Range("H3").Value2 = 47
$defWs.Range("I3").Value2 = 89
$defWs.Range("J3").Value2 = 86
$defWs.Range("N3").Value2 = 28

# ---------------------------------------------------------------------
# Sheet "ST" (special teams) - append this week's kick/FG logs and
# bump the cumulative KO/PT counts
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
$stWs.Range("B2").Value2 = 122
$stWs.Range("D2").Value2 = 120
$stWs.Range("B3").Value2 = 74

$stWs.Range("B6").Value2 = $stWs.Range("B6").Value2 + " 26 22"
$stWs.Range("D3").Value2 = $stWs.Range("D3").Value2 + " 50 47 53 60 41 51"
$stWs.Range("D4").Value2 = $stWs.Range("D4").Value2 + " 7 -1 16 0 0 15"
$stWs.Range("D5").Value2 = $stWs.Range("D5").Value2 + " 0 0 4"

# ---------------------------------------------------------------------
# Sheet "TURNS" - turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value2 = 9
$turnsWs.Range("D2").Value2 = 17
$turnsWs.Range("E2").Value2 = 9

$turnsWs.Range("D3").Value2 = 16
$turnsWs.Range("E3").Value2 = 8

# ---------------------------------------------------------------------
# Sheet "PEN" - penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value2 = 29
$penWs.Range("B3").Value2 = 36
$penWs.Range("B5").Value2 = 1
